$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.104.80'
$ws.Range('E2').Value = '  +0.47%  '
$ws.Range('D3').Value = '1.817.13'
$ws.Range('E3').Value = '  +2.09%  '
$ws.Range('E4').Value = '  -0.32%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '338.10'
$ws.Range('E5').Value = '  -0.80%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.9989'
$ws.Range('E6').Value = '  -0.49%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4320'
$ws.Range('E7').Value = '  +12.95%  '
$ws.Range('E8').Value = '  +2.64%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '45.59'
$ws.Range('E9').Value = '  -3.13%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.153'
$ws.Range('E10').Value = '  +0.83%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07431'
$ws.Range('E11').Value = '  +0.10%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '22.95'
$ws.Range('E12').Value = '  -2.51%  '
$ws.Range('E13').Value = '  -0.21%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.261'
$ws.Range('E14').Value = '  -2.29%  '
$ws.Range('B15').Value = 'WrappedEther'
$ws.Range('C15').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D15').Value = '1.816.07'
$ws.Range('E15').Value = '  +2.04%  '
$ws.Range('B16').Value = 'Chainlink'
$ws.Range('C16').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '7.279'
$ws.Range('E16').Value = '  -1.76%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001084'
$ws.Range('E17').Value = '  +0.52%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.06686'
$ws.Range('E18').Value = '  -0.09%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '81.93'
$ws.Range('E19').Value = '  -0.71%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.9993'
$ws.Range('E20').Value = '  -0.19%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.480'
$ws.Range('E21').Value = '  +0.74%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '17.26'
$ws.Range('E22').Value = '  -1.34%  '
$ws.Range('D23').Value = '28.134.92'
$ws.Range('E23').Value = '  +0.58%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.95'
$ws.Range('E24').Value = '  -1.28%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.385'
$ws.Range('E25').Value = '  -0.68%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.488'
$ws.Range('E26').Value = '  +2.18%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '20.72'
$ws.Range('E27').Value = '  -0.46%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '156.63'
$ws.Range('E28').Value = '  +1.40%  '
$ws.Range('D29').Value = '2.025.24'
$ws.Range('E29').Value = '  +2.22%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.302'
$ws.Range('E30').Value = '  -11.29%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '132.67'
$ws.Range('E31').Value = '  -1.41%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.058'
$ws.Range('E32').Value = '  +0.63%  '
$ws.Range('E33').Value = '  -1.56%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.09231'
$ws.Range('E34').Value = '  +4.37%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '12.36'
$ws.Range('E35').Value = '  -3.06%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.02374'
$ws.Range('E36').Value = '  -1.25%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.6741'
$ws.Range('E37').Value = '  -1.80%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '5.238'
$ws.Range('E38').Value = '  -1.25%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.06255'
$ws.Range('E39').Value = '  -2.37%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.2162'
$ws.Range('E40').Value = '  -0.12%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.483'
$ws.Range('E41').Value = '  -1.20%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.218'
$ws.Range('E42').Value = '  -1.54%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '8.199'
$ws.Range('E43').Value = '  -0.37%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.9980'
$ws.Range('E44').Value = '  -0.52%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '14.02'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.875'
$ws.Range('E46').Value = '  -0.02%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.6128'
$ws.Range('E47').Value = '  -2.28%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '128.95'
$ws.Range('E48').Value = '  -3.59%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.043'
$ws.Range('E49').Value = '  -1.78%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.178'
$ws.Range('E50').Value = '  -2.76%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.07109'
$ws.Range('E51').Value = '  -4.50%  '
